$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value2 = 305.18182
$ws.Range("I38").Value2 = 122.42857
$ws.Range("K38").Value2 = 367.28571
$ws.Range("M38").Value2 = 4.714290000000005

$ws.Range("H70").Value2 = 2116.6667
$ws.Range("I70").Value2 = 1550
$ws.Range("J70").Value2 = 3250
$ws.Range("K70").Value2 = 4650
$ws.Range("L70").Value2 = 9750
$ws.Range("M70").Value2 = -4380
$ws.Range("N70").Value2 = -10290

$ws.Range("H73").Value2 = 2116.6667
$ws.Range("I73").Value2 = 1550
$ws.Range("J73").Value2 = 3250
$ws.Range("K73").Value2 = 4650
$ws.Range("L73").Value2 = 9750
$ws.Range("M73").Value2 = -3714
$ws.Range("N73").Value2 = -11622

$ws.Range("H74").Value2 = 4687.091
$ws.Range("I74").Value2 = 6172.5
$ws.Range("J74").Value2 = 3838.2856
$ws.Range("K74").Value2 = 6172.5
$ws.Range("L74").Value2 = 3838.2856
$ws.Range("M74").Value2 = -5236.5
$ws.Range("N74").Value2 = -5710.2856

$ws.Range("H77").Value2 = 4687.091
$ws.Range("I77").Value2 = 6172.5
$ws.Range("J77").Value2 = 3838.2856
$ws.Range("K77").Value2 = 30862.5
$ws.Range("L77").Value2 = 19191.428
$ws.Range("M77").Value2 = -26182.5
$ws.Range("N77").Value2 = -28551.428

$ws.Range("H88").Value2 = 17545838
$ws.Range("I88").Value2 = 1841.2
$ws.Range("J88").Value2 = 23811552
$ws.Range("K88").Value2 = 1841.2
$ws.Range("L88").Value2 = 23811552
$ws.Range("M88").Value2 = -1435.2
$ws.Range("N88").Value2 = -23812364

$ws.Range("H91").Value2 = 17545838
$ws.Range("I91").Value2 = 1841.2
$ws.Range("J91").Value2 = 23811552
$ws.Range("K91").Value2 = 1841.2
$ws.Range("L91").Value2 = 23811552
$ws.Range("M91").Value2 = -437.2
$ws.Range("N91").Value2 = -23814360

$ws.Range("H130").Value2 = 0
$ws.Range("J130").Value2 = 0
$ws.Range("L130").Value2 = 0
$ws.Range("N130").ClearContents()

$ws.Range("H132").Value2 = 2133749.2
$ws.Range("I132").Value2 = 2233714.2
$ws.Range("J132").Value2 = 1166.6666
$ws.Range("K132").Value2 = 6701142.600000001
$ws.Range("L132").Value2 = 3499.9998
$ws.Range("M132").Value2 = -6698612.600000001
$ws.Range("N132").Value2 = -8559.9998

$ws.Range("H135").Value2 = 27782392
$ws.Range("I135").Value2 = 41667564
$ws.Range("J135").Value2 = 12048.667
$ws.Range("K135").Value2 = 375008076
$ws.Range("L135").Value2 = 108438.003
$ws.Range("M135").Value2 = -375005541
$ws.Range("N135").Value2 = -113508.003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 47030
$ws.Range("I32").Value2 = 45893.1
$ws.Range("K32").Value2 = 45893.1
$ws.Range("M32").Value2 = -45606.1

$ws.Range("H45").Value2 = 1073.5
$ws.Range("I45").Value2 = 985.375
$ws.Range("J45").Value2 = 1249.75
$ws.Range("K45").Value2 = 985.375
$ws.Range("L45").Value2 = 1249.75
$ws.Range("M45").Value2 = -608.375
$ws.Range("N45").Value2 = -2003.75

$ws.Range("H61").Value2 = 2034.8572
$ws.Range("I61").Value2 = 1819.04
$ws.Range("K61").Value2 = 1819.04
$ws.Range("M61").Value2 = -1607.04

$ws.Range("H113").Value2 = 30000
$ws.Range("J113").Value2 = 30000
$ws.Range("L113").Value2 = 30000
$ws.Range("N113").Value2 = -38678

$ws.Range("H136").Value2 = 2034.8572
$ws.Range("I136").Value2 = 1819.04
$ws.Range("K136").Value2 = 5457.12
$ws.Range("M136").Value2 = -2907.12

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value2 = 2196.6
$ws.Range("I10").Value2 = 997.5
$ws.Range("J10").Value2 = 2996
$ws.Range("K10").Value2 = 997.5
$ws.Range("L10").Value2 = 2996
$ws.Range("M10").Value2 = -857.5
$ws.Range("N10").Value2 = -3276

$ws.Range("H88").Value2 = 20114.334
$ws.Range("I88").Value2 = 12000
$ws.Range("J88").Value2 = 24171.5
$ws.Range("K88").Value2 = 12000
$ws.Range("L88").Value2 = 24171.5
$ws.Range("M88").Value2 = -11594
$ws.Range("N88").Value2 = -24983.5

$ws.Range("H91").Value2 = 20114.334
$ws.Range("I91").Value2 = 12000
$ws.Range("J91").Value2 = 24171.5
$ws.Range("K91").Value2 = 12000
$ws.Range("L91").Value2 = 24171.5
$ws.Range("M91").Value2 = -10596
$ws.Range("N91").Value2 = -26979.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 41673404
$ws.Range("I31").Value2 = 7030.478
$ws.Range("K31").Value2 = 7030.478
$ws.Range("M31").Value2 = -6735.478

$ws.Range("H34").Value2 = 41673404
$ws.Range("I34").Value2 = 7030.478
$ws.Range("K34").Value2 = 7030.478
$ws.Range("M34").Value2 = -6828.478

$ws.Range("H58").Value2 = 1329.1555
$ws.Range("I58").Value2 = 1257.1282
$ws.Range("J58").Value2 = 1797.3334
$ws.Range("K58").Value2 = 1257.1282
$ws.Range("L58").Value2 = 1797.3334
$ws.Range("M58").Value2 = -1054.1282
$ws.Range("N58").Value2 = -2203.3334

$ws.Range("H134").Value2 = 2816.9812
$ws.Range("I134").Value2 = 3074.442
$ws.Range("J134").Value2 = 1709.9
$ws.Range("K134").Value2 = 9223.326000000001
$ws.Range("L134").Value2 = 5129.700000000001
$ws.Range("M134").Value2 = -6688.326000000001
$ws.Range("N134").Value2 = -10199.7

$ws.Range("H136").Value2 = 1329.1555
$ws.Range("I136").Value2 = 1257.1282
$ws.Range("J136").Value2 = 1797.3334
$ws.Range("K136").Value2 = 3771.3846
$ws.Range("L136").Value2 = 5392.0002
$ws.Range("M136").Value2 = -1221.3846
$ws.Range("N136").Value2 = -10492.0002

$ws.Range("H141").Value2 = 65589.5
$ws.Range("J141").Value2 = 71733.22
$ws.Range("L141").Value2 = 71733.22
$ws.Range("N141").Value2 = -82093.22

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value2 = 1681.4445
$ws.Range("I34").Value2 = 625
$ws.Range("J34").Value2 = 1983.2858
$ws.Range("K34").Value2 = 1875
$ws.Range("L34").Value2 = 5949.857400000001
$ws.Range("M34").Value2 = -1791
$ws.Range("N34").Value2 = -6117.857400000001

$ws.Range("H92").Value2 = 325
$ws.Range("I92").Value2 = 100
$ws.Range("J92").Value2 = 400
$ws.Range("K92").Value2 = 300
$ws.Range("L92").Value2 = 1200
$ws.Range("M92").Value2 = 948
$ws.Range("N92").Value2 = -3696

$ws.Range("H131").Value2 = 733.92
$ws.Range("J131").Value2 = 804
$ws.Range("L131").Value2 = 2412
$ws.Range("N131").Value2 = -12492

$ws.Range("H133").Value2 = 2237.7666
$ws.Range("I133").Value2 = 2635.5557
$ws.Range("J133").Value2 = 1641.0834
$ws.Range("K133").Value2 = 7906.6671
$ws.Range("L133").Value2 = 4923.2502
$ws.Range("M133").Value2 = -2846.6671
$ws.Range("N133").Value2 = -15043.2502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value2 = 8335049
$ws.Range("I113").Value2 = 62500744
$ws.Range("K113").Value2 = 62500744
$ws.Range("M113").Value2 = -62498574

$ws.Range("H132").Value2 = 6576.8276
$ws.Range("I132").Value2 = 8714.177
$ws.Range("J132").Value2 = 3548.9167
$ws.Range("K132").Value2 = 26142.531
$ws.Range("L132").Value2 = 10646.7501
$ws.Range("M132").Value2 = -23612.531
$ws.Range("N132").Value2 = -15706.7501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 786.875
$ws.Range("I22").Value2 = 877.25
$ws.Range("J22").Value2 = 696.5
$ws.Range("K22").Value2 = 877.25
$ws.Range("L22").Value2 = 696.5
$ws.Range("M22").Value2 = -582.25
$ws.Range("N22").Value2 = -1286.5

$ws.Range("H27").Value2 = 786.875
$ws.Range("I27").Value2 = 877.25
$ws.Range("J27").Value2 = 696.5
$ws.Range("K27").Value2 = 877.25
$ws.Range("L27").Value2 = 696.5
$ws.Range("M27").Value2 = -770.25
$ws.Range("N27").Value2 = -910.5

$ws.Range("H43").Value2 = 6678000
$ws.Range("J43").Value2 = 17000
$ws.Range("L43").Value2 = 17000
$ws.Range("N43").Value2 = -17386

$ws.Range("H46").Value2 = 1351
$ws.Range("I46").Value2 = 850.25
$ws.Range("K46").Value2 = 850.25
$ws.Range("M46").Value2 = -662.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value2 = 933.7692
$ws.Range("I81").Value2 = 739.9091
$ws.Range("J81").Value2 = 2000
$ws.Range("K81").Value2 = 1479.8182
$ws.Range("L81").Value2 = 4000
$ws.Range("M81").Value2 = -418.8181999999999
$ws.Range("N81").Value2 = -6122

$ws.Range("H84").Value2 = 933.7692
$ws.Range("I84").Value2 = 739.9091
$ws.Range("J84").Value2 = 2000
$ws.Range("K84").Value2 = 7399.090999999999
$ws.Range("L84").Value2 = 20000
$ws.Range("M84").Value2 = -2095.090999999999
$ws.Range("N84").Value2 = -30608

$ws.Range("H132").Value2 = 3564.4546
$ws.Range("I132").Value2 = 2958.7144
$ws.Range("J132").Value2 = 4624.5
$ws.Range("K132").Value2 = 8876.143199999999
$ws.Range("L132").Value2 = 13873.5
$ws.Range("M132").Value2 = -6346.143199999999
$ws.Range("N132").Value2 = -18933.5

$ws.Range("H136").Value2 = 2007.2188
$ws.Range("I136").Value2 = 1945.0333
$ws.Range("J136").Value2 = 2940
$ws.Range("K136").Value2 = 2940
$ws.Range("L136").Value2 = 8820
$ws.Range("M136").Value2 = -3285.0999
$ws.Range("N136").Value2 = -13920

